# Code updated for 4 problem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: CN / Implement Lower Bound / Java / Easy
$ws.Range("B3").Value = "CN"
$ws.Range("C3").Value = "Implement Lower Bound"
$ws.Range("D3").Value = "Java"
$ws.Range("E3").Value = "Easy"

# Row 4: CN / Implement Upper Bound / Java / Easy
$ws.Range("B4").Value = "CN"
$ws.Range("C4").Value = "Implement Upper Bound"
$ws.Range("D4").Value = "Java"
$ws.Range("E4").Value = "Easy"

# Row 5: 35 / LC/CN / Search Insert Position / Java / Easy
$ws.Range("A5").Value = 35
$ws.Range("B5").Value = "LC/CN"
$ws.Range("C5").Value = "Search Insert Position"
$ws.Range("D5").Value = "Java"
$ws.Range("E5").Value = "Easy"

# Row 6: CN / Ceil the floor / Java / Easy
$ws.Range("B6").Value = "CN"
$ws.Range("C6").Value = "Ceil the floor"
$ws.Range("D6").Value = "Java"
$ws.Range("E6").Value = "Easy"

# Match the "Easy" level formatting (green fill) already used in E2.
$ws.Range("E2").Copy()
$ws.Range("E3:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to E11, matching the author's final cursor position.
$ws.Range("E11").Select()
